# "case with 380 kV done"
# Refresh the per-bus voltage magnitude (vm_pu) results after re-running the
# power-flow case with the slack/gen bus setpoint changed from 1.05 to 1.02 p.u.
# (column A = bus index, column G = slack bus value which stays 1, column H is
# unused in this export, matching the original layout).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ "B"=1.02; "C"=1.033294655339615; "D"=1.036824891749647; "E"=1.041934602809719; "F"=1.051219653389859; "I"=1.035685900914938; "J"=1.038420169120667; "K"=1.0396175275634; "L"=1.044712713289482; "M"=1.053971761247083; "N"=1.005712725503983 }
    3 = @{ "B"=1.02; "C"=1.034304253336271; "D"=1.037583992375746; "E"=1.04283920727192; "F"=1.052228572741494; "I"=1.035899657571543; "J"=1.039072115948; "K"=1.040186443568828; "L"=1.045427800484763; "M"=1.054792778118291 }
    4 = @{ "B"=1.02; "C"=1.03495767743722; "D"=1.038074878476643; "E"=1.043425045519667; "F"=1.05288193218208; "I"=1.036036196961773; "J"=1.039493533436225; "K"=1.040553607319118; "L"=1.045890376928189; "M"=1.055323940730805 }
    5 = @{ "B"=1.02; "C"=1.035232411664514; "D"=1.038281173467867; "E"=1.043671450427632; "F"=1.053156728420251; "I"=1.036093172543156; "J"=1.03967059240783; "K"=1.040707731653061; "L"=1.046084811579947; "M"=1.05554721913392 }
    6 = @{ "B"=1.02; "C"=1.035278542775973; "D"=1.038315806962615; "E"=1.043712829819853; "F"=1.053202875160731; "I"=1.036102714022394; "J"=1.039700315216076; "K"=1.040733596220351; "L"=1.046117456088337; "M"=1.055584707226854 }
    7 = @{ "B"=1.02; "C"=1.034961348313403; "D"=1.03807763529064; "E"=1.043428337529552; "F"=1.052885603536961; "I"=1.036036959945184; "J"=1.039495899720013; "K"=1.040555667645873; "L"=1.045892975102045; "M"=1.055326924276577 }
    8 = @{ "B"=1.02; "C"=1.033635823107661; "D"=1.037081495424321; "E"=1.042240214345179; "F"=1.051560513997913; "I"=1.035758508281486; "J"=1.038640587439674; "K"=1.039809994311427; "L"=1.044954407432876; "M"=1.054249245961323 }
    9 = @{ "B"=1.02; "C"=1.031301218172145; "D"=1.035323896465381; "E"=1.040150446937716; "F"=1.049229564847008; "I"=1.035254266023179; "J"=1.037130108468288; "K"=1.038488681358278; "L"=1.043299544603147; "M"=1.052349581557568 }
    10 = @{ "B"=1.02; "C"=1.029745600612785; "D"=1.034150695081352; "E"=1.038759911861914; "F"=1.047678353890882; "I"=1.034909007555213; "J"=1.03612093275647; "K"=1.037602914042396; "L"=1.04219567849138; "M"=1.051082739500519 }
    11 = @{ "B"=1.02; "C"=1.029072190603017; "D"=1.033642349577978; "E"=1.038158431636237; "F"=1.047007325114527; "I"=1.034757355115199; "J"=1.03568343605213; "K"=1.037218215397651; "L"=1.041717551926516; "M"=1.050534096663096 }
    12 = @{ "B"=1.02; "C"=1.028822083737021; "D"=1.033453477009692; "E"=1.037935110577587; "F"=1.046758174340672; "I"=1.034700701567712; "J"=1.03552085311438; "K"=1.037075148213442; "L"=1.041539933163563; "M"=1.050330293086309 }
    13 = @{ "B"=1.02; "C"=1.028875731257592; "D"=1.033493993083672; "E"=1.037983009370448; "F"=1.046811613528178; "I"=1.034712868559229; "J"=1.035555731191321; "K"=1.037105844427703; "L"=1.041578033944711; "M"=1.050374010232967 }
    14 = @{ "B"=1.02; "C"=1.029051516127633; "D"=1.033626738335806; "E"=1.038139969895917; "F"=1.046986728206532; "I"=1.034752678703574; "J"=1.035669998467472; "K"=1.037206392944664; "L"=1.041702870334827; "M"=1.05051725046163 }
    15 = @{ "B"=1.02; "C"=1.029159826648912; "D"=1.033708520413284; "E"=1.038236691129626; "F"=1.047094635303011; "I"=1.034777164248876; "J"=1.035740392094691; "K"=1.037268321301869; "L"=1.041779783353209; "M"=1.05060550376309 }
    16 = @{ "B"=1.02; "C"=1.02979029647876; "D"=1.034184425178989; "E"=1.03879984347292; "F"=1.047722901734567; "I"=1.034919026897878; "J"=1.036149957121585; "K"=1.037628420932184; "L"=1.042227407163878; "M"=1.051119149234756 }
    17 = @{ "B"=1.02; "C"=1.030185822557229; "D"=1.034482856857909; "E"=1.039153263417019; "F"=1.048117172715619; "I"=1.035007437403656; "J"=1.036406728380756; "K"=1.037853992988704; "L"=1.042508151258059; "M"=1.051441321300345 }
    18 = @{ "B"=1.02; "C"=1.030416543840744; "D"=1.034656893991253; "E"=1.039359468089573; "F"=1.048347207539535; "I"=1.035058797970868; "J"=1.036556448800317; "K"=1.037985453816939; "L"=1.042671890486432; "M"=1.051629229911371 }
    19 = @{ "B"=1.02; "C"=1.03049521681899; "D"=1.034716230528607; "E"=1.039429788844286; "F"=1.048425654174722; "I"=1.035076275340185; "J"=1.036607491113567; "K"=1.038030259624452; "L"=1.042727718931538; "M"=1.051693300363048 }
    20 = @{ "B"=1.02; "C"=1.030143384529145; "D"=1.03445084137149; "E"=1.039115338507707; "F"=1.048074864625507; "I"=1.03499797328218; "J"=1.036379184413724; "K"=1.037829802773895; "L"=1.04247803152307; "M"=1.051406756203265 }
    21 = @{ "B"=1.02; "C"=1.028999751101738; "D"=1.03358764955211; "E"=1.038093746291779; "F"=1.046935158554891; "I"=1.034740964523495; "J"=1.035636351729405; "K"=1.037176788682553; "L"=1.041666109712067; "M"=1.050475070143973 }
    22 = @{ "B"=1.02; "C"=1.02828086310856; "D"=1.033044635356048; "E"=1.037451983485711; "F"=1.046219154330536; "I"=1.034577503887695; "J"=1.03516885702379; "K"=1.036765212311534; "L"=1.041155499466112; "M"=1.049889206073509 }
    23 = @{ "B"=1.02; "C"=1.028661944043647; "D"=1.033332524762199; "E"=1.037792141452256; "F"=1.046598667066877; "I"=1.034664334485835; "J"=1.035416726956573; "K"=1.036983491359284; "L"=1.04142619507201; "M"=1.050199790711672 }
    24 = @{ "B"=1.02; "C"=1.030162560392776; "D"=1.034465307893755; "E"=1.039132474953913; "F"=1.048093981634902; "I"=1.035002250353314; "J"=1.036391630501003; "K"=1.037840733635334; "L"=1.042491641378226; "M"=1.051422374710289 }
    25 = @{ "B"=1.02; "C"=1.031904632432419; "D"=1.035778541144434; "E"=1.040690239645908; "F"=1.049831688891109; "I"=1.035386230317398; "J"=1.03752099201938; "K"=1.038831137868651; "L"=1.043727479190696; "M"=1.052840763820664 }
}

foreach ($rowKey in $data.Keys) {
    $rowData = $data[$rowKey]
    foreach ($colLetter in $rowData.Keys) {
        $ws.Range($colLetter + $rowKey).Value = $rowData[$colLetter]
    }
}

Write-Output "Updated $($data.Count) rows of vm_pu values"